$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.430.96"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +3.02%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'2.541.85"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.45%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "'  +0.01%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'595.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +2.35%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'177.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +3.60%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("E7").Value = "'  -0.07%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.532"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +1.55%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'2.541.93"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +1.47%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("E10").Value = "'  +2.78%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("E11").Value = "'  +3.00%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("E12").Value = "'  +0.93%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("E13").Value = "'  -0.23%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'27.06"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +1.63%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'3.005.30"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.19%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("E16").Value = "'  +2.82%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'68.382.16"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +3.04%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'2.541.19"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +1.55%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("E19").Value = "'  +4.69%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("E20").Value = "'  +3.02%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'369.02"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +6.27%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("E22").Value = "'  +0.89%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'4.73"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +2.45%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("E24").Value = "'  -0.69%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("B25").Value = "'Dai"
$ws.Range("B25").Style = "Normal"
$ws.Range("C25").Value = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("C25").Style = "Normal"
$ws.Range("D25").Value = "'1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.06%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("B26").Value = "'Aptos"
$ws.Range("B26").Style = "Normal"
$ws.Range("C26").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("C26").Style = "Normal"
$ws.Range("D26").Value = "'10.34"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +3.48%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("E27").Value = "'  +2.09%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'2.677.64"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +1.93%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("E29").Value = "'  -0.17%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("E30").Value = "'  +2.60%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'543.90"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +3.05%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("E32").Value = "'  +2.85%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'1.35"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +2.57%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("E34").Value = "'  +2.95%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("E35").Value = "'  -0.51%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("E36").Value = "'  +0.11%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("E37").Value = "'  +1.52%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'157.65"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.89%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'18.96"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +1.92%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'18.72"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +1.98%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("E41").Value = "'  +2.21%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("B42").Value = "'RenderToken"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'5.25"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +3.08%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("B43").Value = "'PolygonEcosystemToken"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'0.357"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.94%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'2.57"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +2.54%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("E45").Value = "'  +0.10%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("E46").Value = "'  +1.73%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'147.98"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.13%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'0.0₆0282"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +3.36%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("E49").Value = "'  +2.02%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("E51").Value = "'  +0.99%  "
$ws.Range("E51").Style = "Normal"
